$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2Values = @(85.14095539, 88.01228221, 94.3762378, 94.4938712, 89.01497834, 84.44051648999999, 88.54804912, 91.91287499000001, 95.40556424, 99.03097568, 102.7941528, 106.7003306, 110.7549431, 114.963631, 119.3322489, 123.8668744, 128.5738156, 133.4596206, 138.5310862, 143.7952675, 149.2594876, 154.9313482, 160.8187394, 166.9298515, 173.2731859, 179.8575669, 186.6921545, 193.7864563, 201.1503417, 208.7940546, 216.7282287, 224.9639014, 233.5125297, 242.3860058, 251.596674, 261.1573476)
$row3Values = @(6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500, 6552500)
$row10Values = @(3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653, 3.653)
$row11Values = @(17441469.6, 17502271.79, 17553000.5, 17666250.8, 17750912.45, 17817297.53, 17977378.4, 17960842.94, 17968527.09, 17967269.73, 17957849.84, 17938640.19, 17910684.98, 17873450.95, 17828462.18, 17774284.29, 17711215.33, 17639827.5, 17560523.02, 17472114.79, 17374627.32, 17268906.21, 17153965.85, 17030110.77, 16899825.84, 16764039.27, 16623199.05, 16476457.36, 16325052.93, 16167630.76, 16005364.68, 15837808.2, 15665144.12, 15488419.4, 15307295, 15120549.8)
$row12Values = @(3895227.404, 3923222.206, 3953812.502, 4003749.2, 4052087.55, 4101702.47, 4178621.6, 4220157.06, 4272776.911, 4329054.268, 4389180.162, 4453031.811, 4520571.022, 4592487.046, 4668587.816, 4749043.709, 4834184.668, 4923844.504, 5018388.981, 5118021.212, 5222398.679, 5331711.792, 5446184.147, 5565233.227, 5684768.156, 5804672.726, 5924680.949, 6044978.637, 6165273.066, 6284917.236, 6404521.32, 6523449.796, 6641931.879, 6759622.601, 6876868.998, 6992306.196)

$rows = @(2, 3, 10, 11, 12)
$rowData = @{
    2 = $row2Values
    3 = $row3Values
    10 = $row10Values
    11 = $row11Values
    12 = $row12Values
}

foreach ($r in $rows) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 10 + $i  # column J = 10
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}

Write-Host "Updated rows 2, 3, 10, 11, 12 (columns J:AS)"